$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptocurrency price/volume data refresh (GitHub Actions scheduled update)

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '45.164.24'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +4.98%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.362.08'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +2.51%  '

$ws.Range("E4").Value = '  +0.37%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '108.80'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.24%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '309.29'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.57%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.627'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.14%  '

$ws.Range("E8").Value = '  +0.03%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.612'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.29%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '40.81'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.37%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0912'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.03%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '8.46'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.38%  '

$ws.Range("E13").Value = '  +1.49%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.978'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.12%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.719.23'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.55%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '15.36'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.45%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.370.04'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +3.13%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '45.118.42'
$ws.Range("D18").Style = "Normal"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.28'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.79%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0000106'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.92%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.01'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -4.24%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '73.30'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.16%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.43'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.24%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '259.62'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.74%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.27'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.00%  '

$ws.Range("E26").Value = '  -0.55%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.04'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.15%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.39'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -6.27%  '

$ws.Range("E29").Value = '  +4.36%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '22.40'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.44%  '

$ws.Range("B31").Value = 'Hedera'
$ws.Range("C31").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0955'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +10.14%  '

$ws.Range("B32").Value = 'InjectiveProtocol'
$ws.Range("C32").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '37.62'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.82%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '169.11'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.18%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.92'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +3.69%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.130'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.32%  '

$ws.Range("B36").Value = 'Kaspa'
$ws.Range("C36").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.115'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.73%  '

$ws.Range("B37").Value = 'RenderToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.77'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +3.17%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.98'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +6.73%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.91'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +7.91%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0354'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.90%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.72'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +9.58%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '101.26'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -5.04%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.233'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.60%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '13.09'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +6.73%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '69.56'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.88%  '

$ws.Range("E46").Value = '  +0.40%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '81.28'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +6.75%  '

$ws.Range("B48").Value = 'FraxShare'
$ws.Range("C48").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.36'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +5.88%  '

$ws.Range("B49").Value = 'Aave'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '112.43'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.61%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '5.52'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +6.88%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.630.14'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.98%  '
